$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Horas insumidas")
$ws.Activate()

$ws.Cells.Item(72, 2).Value = 40467
$ws.Cells.Item(72, 2).NumberFormat = "d-mmm"
$ws.Cells.Item(72, 3).Value = "Nico"
$ws.Cells.Item(72, 4).Value = "Desarrollo Metricas Agentes"
$ws.Cells.Item(72, 5).Value = "S-01004"
$ws.Cells.Item(72, 6).Value = 1

$ws.Range("F73").Select()
